$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "rxxx"
$ws.Range("B3").Value = "kj"
$ws.Range("C3").Value = "kjhgf"
$ws.Range("D3").Value = "2025-09-27 00:50:48"
